$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column D (uses same style as A3/A5 -> cellXfs idx 5)
$ws.Range("A3").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "Retrained Models Based on Hyperparamter Search V2 "

# New grid-search random forest accuracy values in column D
$ws.Range("C2").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("D2").Value = 0.65

$ws.Range("A6").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("D4").PasteSpecial(-4122)

$ws.Range("C2").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("D5").Value = 0.51

$ws.Range("A6").Copy()
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("D8").PasteSpecial(-4122)

$null = $ws.Range("D3").Select()
